$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.0019918910693377256
$ws.Range("A3").Value = 0.001991890836507082
$ws.Range("G3").Value = 260.0
$ws.Range("H3").Value = 70.0
$ws.Range("I3").Value = 3.7142860889434814
$ws.Range("A4").Value = 0.001991890836507082
